$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (93) with the next day's data, matching the style of
# the existing date column (A) which uses the date-formatted style from A92.
$row = 93

$ws.Cells.Item($row, 1).Value = 46042
$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat

$ws.Cells.Item($row, 2).Value = 216
$ws.Cells.Item($row, 3).Value = 222
$ws.Cells.Item($row, 4).Value = 210
